$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.210.98'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.645.41'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.73'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.60'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.26'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.00'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.127.25'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.278.38'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.649.72'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.38'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '363.64'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.07'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.48'
$ws.Range('E24').Value = '  +3.79%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.75'
$ws.Range('E26').Value = '  -2.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.781.81'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000104'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '555.98'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.07'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.84'
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.40'
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.75'
$ws.Range('E38').Value = '  +2.53%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '158.86'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.72'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.02'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.614'
$ws.Range('E51').Value = '  -0.61%  '
